# "Upd: Fixed incorrect doc"
#
# This fixes a documentation mistake in the SwitchableRangedWeapon sheet:
#   - Row 6 (switchableFiremodes): the "support" (F) and "default" (G)
#     columns incorrectly showed "/"; they should show the actual list of
#     supported fire-mode strings and the actual default value "Auto".
#   - Row 9 (burstReload): the "support" (F) column incorrectly showed the
#     fire-mode string list; it should just show "/". A clarifying remark
#     is also added in the "comment" (H) column explaining that this field
#     only matters in Burst fire mode.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data corrections -----------------------------------------------------

# switchableFiremodes row: fix "support" and "default" columns
$ws.Range("F6").Value2 = "Auto, Semi, Burst, Safe"
$ws.Range("G6").Value2 = "Auto"

# burstReload row: fix "support" column and add explanatory comment
$ws.Range("F9").Value2 = "/"
$ws.Range("H9").Value2 = "此参数只对Burst开火模式有效"

# --- Cosmetic cleanup -------------------------------------------------------
# A number of cells in the table carry a stray/duplicate "no border" cell
# style left over from earlier edits. Re-apply the (identical-looking)
# format that the rest of the table already uses so the workbook's style
# table collapses the duplicate, matching the rest of the sheet.
$formatSource = $ws.Range("B7")
$formatSource.Copy()

$cellsToRestyle = @( `
    "C3", "D3", "E3", "F3", "G3", `
    "C4", "D4", "E4", "F4", "G4", `
    "B6", "C6", "D6", "E6", "F6", "G6", `
    "G7", `
    "G8", `
    "G9", `
    "G10" `
)
foreach ($addr in $cellsToRestyle) {
    $ws.Range($addr).PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = 0

# --- View state -------------------------------------------------------------
$ws.Range("H18").Select()
